# Populate Sheet1 with the subreddit / ideology lookup table used by the
# bot's read_column() helper, replacing the old "Hello"/"Test" placeholder
# content, and make the header row bold.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Subreddit",      "Ideology"),
    @("r/socialism",    "Left"),
    @("r/Libertarian",  "Center Right"),
    @("r/The_Donald",   "Right"),
    @("r/politics",     "Center")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Header row (Subreddit / Ideology) is bold.
$ws.Range("A1:B1").Font.Bold = $true

# Match the resulting sheet view/print setup.
$ws.PageSetup.Orientation = 1 | Out-Null
$ws.Range("H13").Select() | Out-Null
